# Update "2. Data reporter" block (Organization / Contact person / email /
# phone / website) with refreshed contact details, and move the active
# selection to the Organization field, matching the upstream "Add files via
# upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value2  = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value2  = "Kalymbetova Yryskan"
$ws.Range("B8").Value2  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value2  = "(0312) 32 46 55"
$ws.Range("B10").Value2 = "www.stat.gov.kg"

# The indicator title cell (B2) picks up a plain (non Cyrillic-charset) font
# in the refreshed template.
$ws.Range("B2").Font.Name = "Calibri"

# Active cell moves to the Organization field (B7 in the uploaded file).
$ws.Range("B7").Select()
